$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM table (header in row 1, data rows 2-24) is sorted alphabetically by the
# "Comment" column. Two new components need to be inserted in sorted order:
#   1N4731(A) 4.3V | D3 | (blank) | PTH_HDRV02W64P254_1X02 | 1   -> belongs at row 4
#   SMBJ5.0CA      | D5 | (blank) | SMD_DO214AC_SMB        | 1   -> belongs at row 24 (after shift)
# Everything below each insertion point shifts down by one row.
# Work from the bottom of the sheet upwards so source rows are copied before being overwritten.

# --- Step 1: shift old rows 23-24 down to 25-26 (makes room for the SMBJ5.0CA row at 24) ---
$ws.Range("A24:E24").Copy($ws.Range("A26:E26"))
$ws.Range("A23:E23").Copy($ws.Range("A25:E25"))

# --- Step 2: shift old rows 4-22 down to 5-23 (makes room for the 1N4731(A) row at 4) ---
for ($r = 22; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":E" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":E" + ($r + 1))
    $src.Copy($dst)
}

# --- Step 3: write the new row 4 (1N4731(A) 4.3V / D3 / PTH_HDRV02W64P254_1X02) ---
$ws.Cells.Item(4, 1).Value = "1N4731(A) 4.3V"
$ws.Cells.Item(4, 2).Value = "D3"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = "PTH_HDRV02W64P254_1X02"
$ws.Cells.Item(4, 5).Value = 1

# --- Step 4: write the new row 24 (SMBJ5.0CA / D5 / SMD_DO214AC_SMB) ---
$ws.Cells.Item(24, 1).Value = "SMBJ5.0CA"
$ws.Cells.Item(24, 2).Value = "D5"
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = "SMD_DO214AC_SMB"
$ws.Cells.Item(24, 5).Value = 1
